# BDD.xlsx update: restructure "Feuil1" layout, add DOMAINE_ACTION table,
# rename Q2/R2 headers, add filter note in K3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous layout (rows 1-10, cols A-R) before rewriting it,
# since several blocks shift down by one row and gain/lose columns.
$ws.Range("A1:R10").ClearContents()

# --- TABLE USERS ---
$ws.Range("A1").Value = "TABLE USERS"
$ws.Range("A2").Value = "ID_USER"
$ws.Range("B2").Value = "IDENTIFIANT"
$ws.Range("C2").Value = "MDP"
$ws.Range("D2").Value = "NUM_ID"
$ws.Range("E2").Value = "NOM"
$ws.Range("F2").Value = "ADRESSE"
$ws.Range("G2").Value = "SITE_WEB"
$ws.Range("H2").Value = "NOM_REF"
$ws.Range("I2").Value = "PRENOM_REF"
$ws.Range("J2").Value = "FONCTION_REF"
$ws.Range("K2").Value = "TEL_REF"
$ws.Range("L2").Value = "EMAIL_REF"
$ws.Range("M2").Value = "MISSION"
$ws.Range("N2").Value = "ACTIVITE"
$ws.Range("O2").Value = "VALEUR"
$ws.Range("P2").Value = "PROJET"
$ws.Range("Q2").Value = "ACT_ID"
$ws.Range("R2").Value = "VALIDATION"
$ws.Range("K3").Value = "a filtrer"

# --- CONTACT ---
$ws.Range("A4").Value = "CONTACT"
$ws.Range("A5").Value = "ID_CONTACT"
$ws.Range("B5").Value = "ADDR_MAIL"

# --- OFFRE ---
$ws.Range("A7").Value = "OFFRE"
$ws.Range("A8").Value = "ID_OFFRE"
$ws.Range("B8").Value = "TITRE"
$ws.Range("C8").Value = "DEFINITION"
$ws.Range("D8").Value = "TYPE"
$ws.Range("E8").Value = "ACCES"
$ws.Range("F8").Value = "ADRESSE"
$ws.Range("G8").Value = "MOYEN_ACCES"
$ws.Range("H8").Value = "NOM_REF"
$ws.Range("I8").Value = "PRENOM_REF"
$ws.Range("J8").Value = "FONCTION_REF"
$ws.Range("K8").Value = "TEL_REF"

# --- HORAIRE ---
$ws.Range("A10").Value = "HORAIRE"
$ws.Range("A11").Value = "ID_HORAIRE"
$ws.Range("B11").Value = "ID_OFFRE"
$ws.Range("C11").Value = "H_DEBUT"
$ws.Range("D11").Value = "H_FIN"
$ws.Range("E11").Value = "NB_PLACES"
$ws.Range("F11").Value = "DATE"

# --- INSCRIPTION ---
$ws.Range("A13").Value = "INSCRIPTION"
$ws.Range("A14").Value = "ID_INSCRIPTION"
$ws.Range("B14").Value = "ID_OFFRE"
$ws.Range("C14").Value = "ID_HORAIRE"
$ws.Range("D14").Value = "NOM"
$ws.Range("E14").Value = "PRENOM"
$ws.Range("F14").Value = "ADDR_MAIL"
$ws.Range("G14").Value = "TELEPHONE"
$ws.Range("H14").Value = "CONNAISSANCE"

# --- DOMAINE_ACTION (new table) ---
$ws.Range("A16").Value = "DOMAINE_ACTION"
$ws.Range("A17").Value = "ID_DA"
$ws.Range("B17").Value = "ID_USER"
$ws.Range("C17").Value = "ACT1"
$ws.Range("D17").Value = "ACT2"
$ws.Range("E17").Value = "ACT3"
$ws.Range("F17").Value = "ACT4"

# Update selection to match the authored workbook (active cell G17)
$ws.Range("G17").Select()
